$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (C) column date value from 45324 to 45325 for all
#    existing data rows (rows 2-27).
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = 45325
}

# 2. Build the new row 28 by copying the formatting of row 27 (so that the
#    date columns keep their date number format, and the "Artnamn" column
#    keeps its wrap-text style) and then overwrite the values with the new
#    record's data.
$ws.Range("A27:R27").Copy($ws.Range("A28:R28"))

# There is no "Markägare" value for this record, so remove the copied
# content from F28 (column F has no data in any row of this sheet).
$ws.Range("F28").ClearContents()

$ws.Cells.Item(28, 1).Value = "A 4345-2024"   # A28 Beteckning
$ws.Cells.Item(28, 2).Value = 45324           # B28 Datum
$ws.Cells.Item(28, 3).Value = 45325           # C28 Förändrad
$ws.Cells.Item(28, 4).Value = "OKÄNT"         # D28 Län
$ws.Cells.Item(28, 5).Value = "OKÄNT"         # E28 Kommun
$ws.Cells.Item(28, 7).Value = 9.1             # G28 Area (ha)
$ws.Cells.Item(28, 8).Value = 0               # H28 Fridlysta
$ws.Cells.Item(28, 9).Value = 0               # I28 Signalarter
$ws.Cells.Item(28, 10).Value = 0              # J28 NT
$ws.Cells.Item(28, 11).Value = 0              # K28 VU
$ws.Cells.Item(28, 12).Value = 0              # L28 EN
$ws.Cells.Item(28, 13).Value = 0              # M28 CR
$ws.Cells.Item(28, 14).Value = 0              # N28 RE
$ws.Cells.Item(28, 15).Value = 0              # O28 Rödlistade
$ws.Cells.Item(28, 16).Value = 0              # P28 Hotade
$ws.Cells.Item(28, 17).Value = 0              # Q28 Alla arter
# R28 keeps the copied empty/wrap-text formatted cell (Artnamn column).

# 3. Row 27 ends up with an explicit row height after being edited upstream;
#    make that explicit here too so it matches the target workbook.
$ws.Rows.Item(27).RowHeight = 15
